$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.792.03"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.310.54"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.85"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.70"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  +2.34%  "
$ws.Range("D9").Value = "2.308.82"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.50"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.334"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.77"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D15").Value = "2.722.50"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").Value = "58.664.45"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "2.322.74"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.62"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("E20").Value = "  -3.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.66"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.63"
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.07"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.27"
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("D31").Value = "0.0₃0736"
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.28"
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.05"
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "295.09"
$ws.Range("E41").Value = "  -2.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.93"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0959"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0496"
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.557"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.41"
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0212"
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.00"
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.63"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("E51").Value = "  +0.86%  "
